$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1175.6428
$ws.Range("I53").Value = 1079.091
$ws.Range("J53").Value = 1529.6666
$ws.Range("K53").Value = 1079.091
$ws.Range("L53").Value = 1529.6666
$ws.Range("M53").Value = -442.0909999999999
$ws.Range("N53").Value = -2803.6666

$ws.Range("H96").Value = 732.2
$ws.Range("I96").Value = 702.75
$ws.Range("J96").Value = 850
$ws.Range("K96").Value = 2108.25
$ws.Range("L96").Value = 2550
$ws.Range("M96").Value = -735.25
$ws.Range("N96").Value = -5296

$ws.Range("H103").Value = 29412728
$ws.Range("I103").Value = 1590.5
$ws.Range("J103").Value = 55555960
$ws.Range("K103").Value = 4771.5
$ws.Range("L103").Value = 166667880
$ws.Range("M103").Value = -4185.5
$ws.Range("N103").Value = -166669052

$ws.Range("H112").Value = 6192.75
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 6192.75
$ws.Range("K112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("M112").Value = 18578.25
$ws.Range("N112").Value = -20794.25


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1726.6097
$ws.Range("I2").Value = 1048.56
$ws.Range("J2").Value = 2786.0625
$ws.Range("K2").Value = 1048.56
$ws.Range("L2").Value = 2786.0625
$ws.Range("M2").Value = -935.5599999999999
$ws.Range("N2").Value = -3012.0625

$ws.Range("H45").Value = 2435.3462
$ws.Range("I45").Value = 2200.5
$ws.Range("J45").Value = 3218.1667
$ws.Range("K45").Value = 2200.5
$ws.Range("L45").Value = 3218.1667
$ws.Range("M45").Value = -1823.5
$ws.Range("N45").Value = -3972.1667

$ws.Range("H61").Value = 9262.6
$ws.Range("I61").Value = 19499.5
$ws.Range("J61").Value = 2438
$ws.Range("K61").Value = 19499.5
$ws.Range("L61").Value = 2438
$ws.Range("M61").Value = -19287.5
$ws.Range("N61").Value = -2862

$ws.Range("H74").Value = 785.63635
$ws.Range("I74").Value = 712.8
$ws.Range("K74").Value = 712.8
$ws.Range("M74").Value = 161.2

$ws.Range("H77").Value = 785.63635
$ws.Range("I77").Value = 712.8
$ws.Range("K77").Value = 3564
$ws.Range("M77").Value = 804

$ws.Range("H102").Value = 7250
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 7250
$ws.Range("K102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("M102").Value = 7250
$ws.Range("N102").Value = -10494

$ws.Range("H116").Value = 1726.6097
$ws.Range("I116").Value = 1048.56
$ws.Range("J116").Value = 2786.0625
$ws.Range("K116").Value = 1048.56
$ws.Range("L116").Value = 2786.0625
$ws.Range("M116").Value = 1245.44
$ws.Range("N116").Value = -7374.0625

$ws.Range("H122").Value = 5560.963
$ws.Range("I122").Value = 5888.55
$ws.Range("J122").Value = 4625
$ws.Range("K122").Value = 17665.65
$ws.Range("L122").Value = 13875
$ws.Range("M122").Value = -15215.65
$ws.Range("N122").Value = -18775

$ws.Range("H136").Value = 9262.6
$ws.Range("I136").Value = 19499.5
$ws.Range("J136").Value = 2438
$ws.Range("K136").Value = 58498.5
$ws.Range("L136").Value = 7314
$ws.Range("M136").Value = -55948.5
$ws.Range("N136").Value = -12414

$ws.Range("H141").Value = 65630
$ws.Range("J141").Value = 65630
$ws.Range("L141").Value = 65630
$ws.Range("N141").Value = -75990


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1726.6097
$ws.Range("I3").Value = 1048.56
$ws.Range("J3").Value = 2786.0625
$ws.Range("K3").Value = 1048.56
$ws.Range("L3").Value = 2786.0625
$ws.Range("M3").Value = -934.5599999999999
$ws.Range("N3").Value = -3014.0625

$ws.Range("H64").Value = 147.71428
$ws.Range("I64").Value = 87.40000000000001
$ws.Range("J64").Value = 298.5
$ws.Range("K64").Value = 87.40000000000001
$ws.Range("L64").Value = 298.5
$ws.Range("M64").Value = 137.6
$ws.Range("N64").Value = -748.5

$ws.Range("H67").Value = 147.71428
$ws.Range("I67").Value = 87.40000000000001
$ws.Range("J67").Value = 298.5
$ws.Range("K67").Value = 87.40000000000001
$ws.Range("L67").Value = 298.5
$ws.Range("M67").Value = 692.6
$ws.Range("N67").Value = -1858.5

$ws.Range("H105").Value = 1692.8572
$ws.Range("I105").Value = 1690.909
$ws.Range("J105").Value = 1700
$ws.Range("K105").Value = 1690.909
$ws.Range("L105").Value = 1700
$ws.Range("M105").Value = 56.09099999999989
$ws.Range("N105").Value = -5194


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 8647.166999999999
$ws.Range("J50").Value = 8647.166999999999
$ws.Range("L50").Value = 8647.166999999999
$ws.Range("N50").Value = -9897.166999999999

$ws.Range("H51").Value = 9024
$ws.Range("J51").Value = 9024
$ws.Range("L51").Value = 9024
$ws.Range("N51").Value = -10496

$ws.Range("H60").Value = 7189.5557
$ws.Range("J60").Value = 7952.1665
$ws.Range("L60").Value = 7952.1665
$ws.Range("N60").Value = -8974.166499999999

$ws.Range("H61").Value = 9024
$ws.Range("J61").Value = 9024
$ws.Range("L61").Value = 9024
$ws.Range("N61").Value = -9720

$ws.Range("H68").Value = 16822.428
$ws.Range("J68").Value = 16822.428
$ws.Range("L68").Value = 16822.428
$ws.Range("N68").Value = -18320.428

$ws.Range("H71").Value = 16822.428
$ws.Range("J71").Value = 16822.428
$ws.Range("L71").Value = 50467.284
$ws.Range("N71").Value = -57955.284

$ws.Range("H74").Value = 16375.1
$ws.Range("J74").Value = 16375.1
$ws.Range("L74").Value = 16375.1
$ws.Range("N74").Value = -18123.1

$ws.Range("H77").Value = 16375.1
$ws.Range("J77").Value = 16375.1
$ws.Range("L77").Value = 49125.3
$ws.Range("N77").Value = -57861.3

$ws.Range("H134").Value = 851.5454999999999
$ws.Range("I134").Value = 807.8823
$ws.Range("K134").Value = 2423.6469
$ws.Range("M134").Value = 111.3531000000003

$ws.Range("H138").Value = 49950
$ws.Range("J138").Value = 49950
$ws.Range("L138").Value = 49950
$ws.Range("N138").Value = -60230


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 746.6
$ws.Range("J52").Value = 746.6
$ws.Range("L52").Value = 2239.8
$ws.Range("N52").Value = -2771.8

$ws.Range("H107").Value = 529.86957
$ws.Range("I107").Value = 281.35715
$ws.Range("J107").Value = 916.44446
$ws.Range("K107").Value = 844.0714499999999
$ws.Range("L107").Value = 2749.33338
$ws.Range("M107").Value = 1075.92855
$ws.Range("N107").Value = -6589.33338

$ws.Range("H129").Value = 22223466
$ws.Range("I129").Value = 1407
$ws.Range("J129").Value = 30304216
$ws.Range("K129").Value = 4221
$ws.Range("L129").Value = 90912648
$ws.Range("M129").Value = 779
$ws.Range("N129").Value = -90922648

$ws.Range("H137").Value = 5624.3125
$ws.Range("I137").Value = 5918.048
$ws.Range("J137").Value = 5063.5454
$ws.Range("K137").Value = 17754.144
$ws.Range("L137").Value = 15190.6362
$ws.Range("M137").Value = -12654.144
$ws.Range("N137").Value = -25390.6362


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 18750.092
$ws.Range("J93").Value = 20025.1
$ws.Range("L93").Value = 20025.1
$ws.Range("N93").Value = -23769.1

$ws.Range("H113").Value = 6876895.5
$ws.Range("I113").Value = 20001574
$ws.Range("J113").Value = 911132.9399999999
$ws.Range("K113").Value = 20001574
$ws.Range("L113").Value = 911132.9399999999
$ws.Range("M113").Value = -19999404
$ws.Range("N113").Value = -915472.9399999999

$ws.Range("H140").Value = 89789
$ws.Range("J140").Value = 89789
$ws.Range("L140").Value = 89789
$ws.Range("N140").Value = -100149


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2665.375
$ws.Range("I122").Value = 2493.3044
$ws.Range("J122").Value = 3105.111
$ws.Range("K122").Value = 7479.9132
$ws.Range("L122").Value = 9315.332999999999
$ws.Range("M122").Value = -5029.9132
$ws.Range("N122").Value = -14215.333

$ws.Range("H138").Value = 58987
$ws.Range("J138").Value = 58987
$ws.Range("L138").Value = 58987
$ws.Range("N138").Value = -69267


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 2004088.8
$ws.Range("I29").Value = 10000400
$ws.Range("K29").Value = 10000400
$ws.Range("M29").Value = -10000110

$ws.Range("H138").Value = 78550
$ws.Range("J138").Value = 93066.664
$ws.Range("L138").Value = 93066.664
$ws.Range("N138").Value = -103346.664

